$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Done " status for the "statefullnes of the pages" row (row 6)
$ws.Range("C6").Value = "Done "

# Add new plan item row
$ws.Range("B8").Value = "Klib library implementation "

# Update the active selection to reflect the last edited cell
$ws.Range("C7").Select()
